$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 15 data: GFG | Array Subset of another array | Java | 2023-02-28 (serial 44985)
$ws.Range("A15").Value = "GFG"
$ws.Range("B15").Value = "Array Subset of another array"
$ws.Range("C15").Value = "Java"
$ws.Range("D15").Value = 44985
$ws.Range("D15").NumberFormat = "d-mmm-yy"

# Match formatting of neighboring rows (row 13/14 style for column A/B, row 11-14 style for D)
$ws.Range("A15").HorizontalAlignment = -4108  # xlCenter
$ws.Range("B15").WrapText = $true

# Update selection to E15 as in the diff
$ws.Range("E15").Select()
